$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 367 (shifts existing rows 367:386 down to 368:387,
# carrying over formatting/styles from the surrounding rows).
$ws.Rows(367).Insert()

# Populate the newly inserted row with the new price-report record.
$ws.Cells.Item(367, 1).Value = 4
$ws.Cells.Item(367, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(367, 3).Value = "Los Lagos"
$ws.Cells.Item(367, 4).Value = 44939
$ws.Cells.Item(367, 5).Value = 10
$ws.Cells.Item(367, 6).Value = 100112040
$ws.Cells.Item(367, 7).Value = "Cilantro"
$ws.Cells.Item(367, 8).Value = "Sin especificar"
$ws.Cells.Item(367, 9).Value = "Primera"
$ws.Cells.Item(367, 10).Value = 160
$ws.Cells.Item(367, 11).Value = 8000
$ws.Cells.Item(367, 12).Value = 10000
$ws.Cells.Item(367, 13).Value = 9000
$ws.Cells.Item(367, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(367, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(367, 16).Value = 4500
$ws.Cells.Item(367, 17).Value = 2
$ws.Cells.Item(367, 18).Value = "Hortaliza"
